$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'" + '2026-02-16 06:18:46'
$ws.Range("I2").Value = "'" + '2.4 mm'
$ws.Range("M2").Value = "'" + '2.0 °C 5:48 TU'
$ws.Range("E3").Value = "'" + '2026-02-16 06:18:48'
$ws.Range("I3").Value = "'" + '0.6 mm'
$ws.Range("M3").Value = "'" + '-1.1 °C 5:59 TU'
$ws.Range("O3").Value = "'" + '-1.3 °C'
$ws.Range("E4").Value = "'" + '2026-02-16 06:18:50'
$ws.Range("J4").Value = "'" + '1014.3 hPa'
$ws.Range("O4").Value = "'" + '10.4 °C'
$ws.Range("E5").Value = "'" + '2026-02-16 06:18:53'
$ws.Range("H5").Value = "'" + '97%'
$ws.Range("I5").Value = "'" + '2.7 mm'
$ws.Range("M5").Value = "'" + '-0.6 °C 5:59 TU'
$ws.Range("E6").Value = "'" + '2026-02-16 06:18:55'
$ws.Range("J6").Value = "'" + '1014.4 hPa'
$ws.Range("E7").Value = "'" + '2026-02-16 06:18:58'
$ws.Range("M7").Value = "'" + '14.0 °C 5:37 TU'
$ws.Range("O7").Value = "'" + '13.3 °C'
$ws.Range("E8").Value = "'" + '2026-02-16 06:19:00'
$ws.Range("E9").Value = "'" + '2026-02-16 06:19:03'
$ws.Range("H9").Value = "'" + '97%'
$ws.Range("O9").Value = "'" + '4.9 °C'
$ws.Range("E10").Value = "'" + '2026-02-16 06:19:05'
$ws.Range("K10").Value = "'" + '-0.1 MJ/m2'
$ws.Range("E11").Value = "'" + '2026-02-16 06:19:08'
$ws.Range("N11").Value = "'" + '-0.2 °C 5:33 TU'
$ws.Range("E12").Value = "'" + '2026-02-16 06:19:10'
$ws.Range("O12").Value = "'" + '5.3 °C'
$ws.Range("E13").Value = "'" + '2026-02-16 06:19:12'
$ws.Range("H13").Value = "'" + '89%'
$ws.Range("O13").Value = "'" + '0.7 °C'
$ws.Range("E14").Value = "'" + '2026-02-16 06:19:15'
$ws.Range("E15").Value = "'" + '2026-02-16 06:19:17'
$ws.Range("H15").Value = "'" + '91%'
$ws.Range("N15").Value = "'" + '3.0 °C 5:39 TU'
$ws.Range("O15").Value = "'" + '4.9 °C'
$ws.Range("E16").Value = "'" + '2026-02-16 06:19:20'
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = ""
$ws.Range("I16").Value = ""
$ws.Range("K16").Value = ""
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = ""
$ws.Range("O16").Value = ""
$ws.Range("E17").Value = "'" + '2026-02-16 06:19:36'
$ws.Range("H17").Value = "'" + '59%'
$ws.Range("K17").Value = "'" + '-0.1 MJ/m2'
$ws.Range("L17").Value = "'" + '42.1 km/h - 297º 5:48 TU'
$ws.Range("E18").Value = "'" + '2026-02-16 06:19:42'
$ws.Range("J18").Value = "'" + '1014.8 hPa'
$ws.Range("N18").Value = "'" + '2.3 °C 5:56 TU'
$ws.Range("O18").Value = "'" + '4.0 °C'
$ws.Range("E19").Value = "'" + '2026-02-16 06:19:45'
$ws.Range("N19").Value = "'" + '2.2 °C 5:56 TU'
$ws.Range("O19").Value = "'" + '3.1 °C'
$ws.Range("E20").Value = "'" + '2026-02-16 06:19:47'
$ws.Range("H20").Value = "'" + '93%'
$ws.Range("E21").Value = "'" + '2026-02-16 06:19:50'
$ws.Range("N21").Value = "'" + '3.2 °C 5:59 TU'
$ws.Range("O21").Value = "'" + '4.5 °C'
$ws.Range("E22").Value = "'" + '2026-02-16 06:19:52'
$ws.Range("E23").Value = "'" + '2026-02-16 06:19:55'
$ws.Range("I23").Value = "'" + '1.2 mm'
$ws.Range("L23").Value = "'" + '61.2 km/h - 322º 5:43 TU'
$ws.Range("E24").Value = "'" + '2026-02-16 06:19:57'
$ws.Range("J24").Value = "'" + '1017.9 hPa'
$ws.Range("E25").Value = "'" + '2026-02-16 06:20:00'
$ws.Range("H25").Value = "'" + '73%'
$ws.Range("I25").Value = "'" + '0.1 mm'
$ws.Range("O25").Value = "'" + '0.4 °C'
$ws.Range("E26").Value = "'" + '2026-02-16 06:20:02'
$ws.Range("E27").Value = "'" + '2026-02-16 06:20:04'
$ws.Range("H27").Value = "'" + '76%'
$ws.Range("E28").Value = "'" + '2026-02-16 06:20:07'
$ws.Range("J28").Value = "'" + '1015.8 hPa'
$ws.Range("N28").Value = "'" + '1.3 °C 5:42 TU'
$ws.Range("O28").Value = "'" + '3.0 °C'
$ws.Range("E29").Value = "'" + '2026-02-16 06:20:10'
$ws.Range("K29").Value = "'" + '-0.1 MJ/m2'
$ws.Range("L29").Value = "'" + '7.2 km/h - 0º 5:48 TU'
$ws.Range("M29").Value = "'" + '6.1 °C 5:55 TU'
$ws.Range("E30").Value = "'" + '2026-02-16 06:20:12'
$ws.Range("H30").Value = "'" + '87%'
$ws.Range("J30").Value = "'" + '1014.4 hPa'
$ws.Range("E31").Value = "'" + '2026-02-16 06:20:14'
$ws.Range("J31").Value = "'" + '1013.0 hPa'
$ws.Range("O31").Value = "'" + '13.8 °C'
$ws.Range("E32").Value = "'" + '2026-02-16 06:20:17'
$ws.Range("H32").Value = "'" + '83%'
$ws.Range("E33").Value = "'" + '2026-02-16 06:20:19'
$ws.Range("H33").Value = "'" + '75%'
$ws.Range("O33").Value = "'" + '4.1 °C'
$ws.Range("E34").Value = "'" + '2026-02-16 06:20:22'
$ws.Range("H34").Value = "'" + '65%'
$ws.Range("E35").Value = "'" + '2026-02-16 06:20:24'
$ws.Range("J35").Value = "'" + '1018.9 hPa'
$ws.Range("O35").Value = "'" + '6.8 °C'
$ws.Range("E36").Value = "'" + '2026-02-16 06:20:27'
$ws.Range("E37").Value = "'" + '2026-02-16 06:20:29'
$ws.Range("J37").Value = "'" + '1018.1 hPa'
$ws.Range("E38").Value = "'" + '2026-02-16 06:20:32'
$ws.Range("H38").Value = "'" + '94%'
$ws.Range("L38").Value = "'" + '14.8 km/h - 55º 5:51 TU'
$ws.Range("E39").Value = "'" + '2026-02-16 06:20:34'
$ws.Range("E40").Value = "'" + '2026-02-16 06:20:37'
$ws.Range("O40").Value = "'" + '2.7 °C'
$ws.Range("E41").Value = "'" + '2026-02-16 06:20:39'
$ws.Range("K41").Value = "'" + '-0.1 MJ/m2'
$ws.Range("E42").Value = "'" + '2026-02-16 06:20:41'
$ws.Range("H42").Value = "'" + '96%'
$ws.Range("N42").Value = "'" + '4.8 °C 5:48 TU'
$ws.Range("O42").Value = "'" + '6.3 °C'
$ws.Range("E43").Value = "'" + '2026-02-16 06:20:44'
$ws.Range("N43").Value = "'" + '2.2 °C 5:59 TU'
$ws.Range("O43").Value = "'" + '3.2 °C'
$ws.Range("E44").Value = "'" + '2026-02-16 06:20:46'
$ws.Range("I44").Value = "'" + '1.4 mm'
$ws.Range("E45").Value = "'" + '2026-02-16 06:20:48'
$ws.Range("I45").Value = "'" + '1.4 mm'
$ws.Range("J45").Value = "'" + '1019.4 hPa'
$ws.Range("E46").Value = "'" + '2026-02-16 06:20:51'
$ws.Range("J46").Value = "'" + '1018.4 hPa'
$ws.Range("M46").Value = "'" + '13.4 °C 5:32 TU'
$ws.Range("O46").Value = "'" + '12.7 °C'
